# Weekly Fruta/Hortaliza update: a new price-report row is inserted right
# after the anchor row (434), pushing the existing rows 435:461 down to
# 436:462. The new row 435 carries a fresh "Segunda" quality record dated
# 2022-08-10 (serial 44783) for Acelga @ Terminal La Palmera de La Serena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 435:461 down to 436:462, duplicating row 435's formatting
# (and therefore the D-column date style) into the freshly inserted row.
$ws.Rows(435).Insert()

# Populate the newly inserted row 435 with the new weekly record.
$ws.Range("A435").Value = 8
$ws.Range("B435").Value = "Terminal La Palmera de La Serena"
$ws.Range("C435").Value = "Coquimbo"
$ws.Range("D435").Value = 44783
$ws.Range("E435").Value = 4
$ws.Range("F435").Value = 100112009
$ws.Range("G435").Value = "Acelga"
$ws.Range("H435").Value = "Sin especificar"
$ws.Range("I435").Value = "Segunda"
$ws.Range("J435").Value = 1400
$ws.Range("K435").Value = 600
$ws.Range("L435").Value = 650
$ws.Range("M435").Value = 625
$ws.Range("N435").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O435").Value = "Provincia del Elquí"
$ws.Range("P435").Value = 312
$ws.Range("Q435").Value = 2
$ws.Range("R435").Value = "Hortaliza"
